# normalization correction over dates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.141199946403503
$ws.Range("C2").Value = 0.0009
$ws.Range("E2").Value = 0.3989000022411346

$ws.Range("B3").Value = 0.4153999984264374
$ws.Range("C3").Value = 0.0003
$ws.Range("E3").Value = 0.06069999933242798

$ws.Range("B5").Value = 53.41999816894531
$ws.Range("C5").Value = 0.0564
$ws.Range("E5").Value = 0.7074999809265137

$ws.Range("B6").Value = 113.3384017944336
$ws.Range("C6").Value = 0.0897
$ws.Range("E6").Value = 1

$ws.Range("B7").Value = 22.33880043029785
$ws.Range("C7").Value = 0.0225
$ws.Range("E7").Value = 0.5268999934196472

$ws.Range("B8").Value = 35.45539855957031
$ws.Range("C8").Value = 0.0318
$ws.Range("E8").Value = 0.5268999934196472

$ws.Range("B9").Value = 226.1092071533203
$ws.Range("C9").Value = 0.029
$ws.Range("E9").Value = 1
